$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.108.80"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.873.54"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'313.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'0.5139"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").Value = "'0.3894"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.74%  "

$ws.Range("D9").Value = "'0.08384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").Value = "'1.115"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").Value = "'41.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "'6.200"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").Value = "'20.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "1.867.00"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").Value = "'7.291"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").Value = "'90.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "'0.06650"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "'17.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'6.038"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("D23").Value = "28.150.10"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").Value = "'2.253"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").Value = "2.083.40"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").Value = "'2.474"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.21%  "

$ws.Range("D28").Value = "'158.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").Value = "'20.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").Value = "'125.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").Value = "'0.1061"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("D32").Value = "'1.039"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("D33").Value = "'5.905"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.90%  "

$ws.Range("D34").Value = "'3.594"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").Value = "'9.708"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "'0.02441"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").Value = "'0.06543"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "'0.2186"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("E39").Value = "  -1.23%  "

$ws.Range("D40").Value = "'0.6503"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "'0.6099"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").Value = "'13.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").Value = "'1.279"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").Value = "'3.674"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").Value = "'2.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").Value = "'121.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "'77.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "
